$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Id_habitacion (column A) and id_servicios (column B) values
# by prefixing each existing value with "10004" (new hotel id prefix)
for ($r = 2; $r -le 11; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    $b = $ws.Cells.Item($r, 2).Value2
    $aStr = [string]::Concat("10004", [string]([int]$a))
    $bStr = [string]::Concat("10004", [string]([int]$b))
    $ws.Cells.Item($r, 1).Value = [double]$aStr
    $ws.Cells.Item($r, 2).Value = [double]$bStr
}

# Set column B width to best-fit sized width matching new (longer) values
$ws.Columns.Item(2).ColumnWidth = 10.25

# Update the active selection to C8
$ws.Range("C8").Select()
